$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.115.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.41%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.836.82'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.40%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.28%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.79%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6299'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.10%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.29%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07505'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.23%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2930'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.59%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.23'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.36%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07735'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.15%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.833.97'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.31%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.999'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.27%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6698'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.97%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.72'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.32%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009292'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -7.51%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.007'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.22%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.124.89'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.48%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.084.69'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.35%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.50%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '223.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.84%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.003'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.45%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.144'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.01%  '

# Row 24
$ws.Range("E24").Value = '  +0.49%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.70'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.05%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1404'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.52%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.513'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.70%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.50%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.501'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.95%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05804'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +11.86%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.162'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.98%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.067'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.22%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.206'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.15%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7512'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.89%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.850'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.33%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.140'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.31%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.675'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.75%  '

# Row 38
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.769'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.38%  '

# Row 39
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.226.21'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.27%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01787'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.22%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.564'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.97%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8940'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.05%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.003'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.33%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.85%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.981.79'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.07%  '

# Row 46
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.89'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.23%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000123'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.35%  '

# Row 48
$ws.Range("B48").Value = 'XinFinNetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.07781'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +13.72%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5100'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.12%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4076'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.36%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.067'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.17%  '
